# Add three new data rows (166-168) to the "Espárragos" sheet, mirroring
# the existing rows' structure (Banquete / Primera / Segunda grades for
# the same market date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=166; A=6; B="Mercado Mayorista Lo Valledor de Santiago"; C="Metropolitana"; D=45191; E=13; F=300000000; G="Espárragos"; H="Sin especificar"; I="Banquete"; J=2200; K=1500; L=1600; M=1545; N="`$/kilo"; O="Provincia de Linares"; P=1545; Q=1; R="Hortaliza" },
    @{ Row=167; A=6; B="Mercado Mayorista Lo Valledor de Santiago"; C="Metropolitana"; D=45191; E=13; F=300000000; G="Espárragos"; H="Sin especificar"; I="Primera";  J=1400; K=1200; L=1300; M=1257; N="`$/kilo"; O="Provincia de Linares"; P=1257; Q=1; R="Hortaliza" },
    @{ Row=168; A=6; B="Mercado Mayorista Lo Valledor de Santiago"; C="Metropolitana"; D=45191; E=13; F=300000000; G="Espárragos"; H="Sin especificar"; I="Segunda";  J=1100; K=900;  L=1000; M=955;  N="`$/kilo"; O="Provincia de Linares"; P=955;  Q=1; R="Hortaliza" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C

    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = $r.H
    $ws.Cells.Item($rowNum, 9).Value = $r.I
    $ws.Cells.Item($rowNum, 10).Value = $r.J
    $ws.Cells.Item($rowNum, 11).Value = $r.K
    $ws.Cells.Item($rowNum, 12).Value = $r.L
    $ws.Cells.Item($rowNum, 13).Value = $r.M
    $ws.Cells.Item($rowNum, 14).Value = $r.N
    $ws.Cells.Item($rowNum, 15).Value = $r.O
    $ws.Cells.Item($rowNum, 16).Value = $r.P
    $ws.Cells.Item($rowNum, 17).Value = $r.Q
    $ws.Cells.Item($rowNum, 18).Value = $r.R
}
